# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" quarter sheet right after "总计" and before
# "2022-Q2" (all later quarter sheets shift right by one position), and
# records the corresponding summary row on the "总计" sheet (plus appends
# the 2020-Q4 totals row that the summary table was missing).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    under the header, and append a new row for 2020-Q4 at the bottom.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Insert a fresh, blank row directly below the header (current row 2),
# pushing the existing quarters down by one.
$total.Rows.Item(2).Insert()

# Column A carries a bordered/bold style on every data row; clone it from
# an existing styled cell (now shifted to row 3) instead of re-building
# the font/border/alignment by hand.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.03

# Re-number the index column (A) for the rows that followed, since they
# all moved down one row and the sequence must stay 0,1,2,3,4,5.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Append the brand-new 2020-Q4 row at the end, cloning A7's style for A8.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial($xlPasteFormats)

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2020-Q4"
$total.Range("C8").Value = 2
$total.Range("D8").Value = 0.02

# ---------------------------------------------------------------------
# 2) Insert a brand-new worksheet named "2022-Q3" right before "2022-Q2".
#    Cloning the "2022-Q2" sheet (instead of Worksheets.Add + paste)
#    guarantees the header/index-column styling matches exactly, then we
#    drop the extra rows and overwrite the data with the Q3 figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# "2022-Q2" has 6 data rows (rows 2-7); Q3 only needs 2 (rows 2-3).
$q3.Rows("4:7").Delete()

function Set-FundRow($sheet, $row, $idx, $code, $fundName, $scale, $position, $ratio, $marketValue, $rank) {
    $sheet.Cells.Item($row, 1).Value = $idx
    $sheet.Cells.Item($row, 3).Value = $fundName
    $sheet.Cells.Item($row, 8).Value = $rank

    # Columns B, D-G hold plain text in the source data (e.g. "006923",
    # "0.26"), not numbers, so force a text number-format before
    # assigning the value and reset the style afterwards so no stray
    # number-format sticks around once the text is in place.
    $rng = $sheet.Range($sheet.Cells.Item($row, 2), $sheet.Cells.Item($row, 7))
    $rng.NumberFormat = "@"
    $sheet.Cells.Item($row, 2).Value = $code
    $sheet.Cells.Item($row, 4).Value = $scale
    $sheet.Cells.Item($row, 5).Value = $position
    $sheet.Cells.Item($row, 6).Value = $ratio
    $sheet.Cells.Item($row, 7).Value = $marketValue
    $rng.Style = "Normal"
}

Set-FundRow $q3 2 0 "006923" "前海开源沪港深非周期性行业股票A" "0.26" "86.78" "5.56" "0.0145" 8
Set-FundRow $q3 3 1 "006924" "前海开源沪港深非周期性行业股票C" "0.20" "86.78" "5.56" "0.0111" 8

$total.Activate()
